$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- F column (market_value): replace static numbers with formulas ---
# Row 2 gets its own (non-shared) formula.
$ws.Range("F2").Formula = "=ROUND(E2 * 1.2,0)"
# Rows 3:21 become one shared-formula group.
$ws.Range("F3:F21").Formula = "=ROUND(E3 * 1.2,0)"

# --- J column (carry_max): mirror the I column (amount) values for every item row ---
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 9).Value2()
}

# --- Sheet view: reset scroll position and move the active selection to G10 ---
$ws.Range("G10").Select()
